$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clientes")

# Insert a new row at row 20, pushing existing rows 20+ down by one.
$ws.Rows("20").Insert()

# Populate the new row 20 with the new client record.
$ws.Range("A20").Value = "FAST HORIZON LOGISTICA INTELIGENTE LTDA"
$ws.Range("B20").Value = "SABRINA"
